# Generate Report for Handback
# - Updates the "Status" text from "Ready for handoff" to
#   "Handed back: in sync with en-US" everywhere it appears (Overview E/F
#   columns and the per-language "Status" column).
# - Fills in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns on the zh-cn and de-de sheets with
#   hyperlinked handback file names and a handback timestamp.
# - Widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is used by Overview!E2:F3 as well as the
#    per-language Status column (C2:C3) on zh-cn/de-de.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# 2. zh-cn sheet: row 2 (64c53e2e...) and row 3 (73e28401...)
#    Latest Target File (I) gets a hyperlink to the handback .md file,
#    Latest Handback File (J) gets the generated handback xliff name,
#    Latest Handback DateTime (K) gets the new handback timestamp.
# ---------------------------------------------------------------------
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b02ae85d0352b4815ad8b707d677cda5611c8474/e2e/"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($ghBase + "64c53e2e-e6e1-4ed2-836d-c85d674d0b56.md"), "", "", "64c53e2e-e6e1-4ed2-836d-c85d674d0b56.md")
$zhcn.Range("J2").Value = "64c53e2e-e6e1-4ed2-836d-c85d674d0b56.c7d8faae26e1c2799ee6eb32fe8ad11136b4fb5d.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-16 02:43:12"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($ghBase + "73e28401-ac3f-4dda-8550-b29fa5410a52.md"), "", "", "73e28401-ac3f-4dda-8550-b29fa5410a52.md")
$zhcn.Range("J3").Value = "73e28401-ac3f-4dda-8550-b29fa5410a52.e35c21e5179eb3251e988343801b456cc2c9f908.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-16 02:43:12"

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape, different handback datetime/file names.
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), ($ghBase + "64c53e2e-e6e1-4ed2-836d-c85d674d0b56.md"), "", "", "64c53e2e-e6e1-4ed2-836d-c85d674d0b56.md")
$dede.Range("J2").Value = "64c53e2e-e6e1-4ed2-836d-c85d674d0b56.c7d8faae26e1c2799ee6eb32fe8ad11136b4fb5d.de-de.xlf"
$dede.Range("K2").Value = "2016-08-16 02:43:19"

$dede.Hyperlinks.Add($dede.Range("I3"), ($ghBase + "73e28401-ac3f-4dda-8550-b29fa5410a52.md"), "", "", "73e28401-ac3f-4dda-8550-b29fa5410a52.md")
$dede.Range("J3").Value = "73e28401-ac3f-4dda-8550-b29fa5410a52.e35c21e5179eb3251e988343801b456cc2c9f908.de-de.xlf"
$dede.Range("K3").Value = "2016-08-16 02:43:19"

# ---------------------------------------------------------------------
# 4. Column widths: widen the Status-like column and the newly
#    populated Target/Handback file columns so the longer text fits.
#    (ColumnWidth is specified in characters; the stored OOXML width
#    ends up ~0.833 wider, matching Excel's own MDW padding.)
# ---------------------------------------------------------------------
$wideStatus = 29.144371396019366   # -> stored width ~29.98
$wideFile   = 39.166666666666664   # -> stored width 40

$overview.Columns.Item(5).ColumnWidth = $wideStatus   # E
$overview.Columns.Item(6).ColumnWidth = $wideStatus   # F

$zhcn.Columns.Item(3).ColumnWidth = $wideStatus        # C
$zhcn.Columns.Item(9).ColumnWidth = $wideFile           # I
$zhcn.Columns.Item(10).ColumnWidth = $wideFile          # J

$dede.Columns.Item(3).ColumnWidth = $wideStatus         # C
$dede.Columns.Item(9).ColumnWidth = $wideFile            # I
$dede.Columns.Item(10).ColumnWidth = $wideFile           # J

Write-Host "Handback report generated."
